# Refresh the cryptocurrency price/volume snapshot (row 2..51) with the
# latest scrape results. Prices in column D and 1h volume deltas in column E
# are stored as literal text (not numbers), matching the page markup, so
# numeric-looking price strings are written with a leading apostrophe to
# force text entry, then the quote-prefix formatting flag is cleared again
# by resetting the cell style to "Normal" (keeps the underlying value as
# text without leaving a stray number-format behind).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'30.221.38"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +3.52%  "

$c = $ws.Range("D3")
$c.Value = "'1.906.96"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("E4").Value = "  -0.34%  "

$c = $ws.Range("D5")
$c.Value = "'326.52"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.53%  "

$ws.Range("E6").Value = "  -0.30%  "

$c = $ws.Range("D7")
$c.Value = "'0.5154"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.73%  "

$c = $ws.Range("D8")
$c.Value = "'0.4022"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.36%  "

$c = $ws.Range("D9")
$c.Value = "'0.08467"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.52%  "

$c = $ws.Range("D10")
$c.Value = "'42.63"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.28%  "

$c = $ws.Range("D11")
$c.Value = "'1.121"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.06%  "

$c = $ws.Range("D12")
$c.Value = "'23.23"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +13.14%  "

$c = $ws.Range("D13")
$c.Value = "'6.475"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.81%  "

$c = $ws.Range("D14")
$c.Value = "'1.907.41"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.44%  "

$c = $ws.Range("D15")
$c.Value = "'7.363"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "

$ws.Range("E16").Value = "  -0.37%  "

$c = $ws.Range("D17")
$c.Value = "'94.89"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.98%  "

$c = $ws.Range("D18")
$c.Value = "'0.00001114"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.43%  "

$c = $ws.Range("D19")
$c.Value = "'0.06670"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.85%  "

$c = $ws.Range("D20")
$c.Value = "'18.36"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.64%  "

$ws.Range("E21").Value = "  -0.25%  "

$c = $ws.Range("D22")
$c.Value = "'6.001"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.55%  "

$c = $ws.Range("D23")
$c.Value = "'30.226.31"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +3.44%  "

$c = $ws.Range("D24")
$c.Value = "'11.29"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.04%  "

$c = $ws.Range("D25")
$c.Value = "'2.202"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.70%  "

$c = $ws.Range("D26")
$c.Value = "'2.125.14"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "

$c = $ws.Range("D27")
$c.Value = "'21.61"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.27%  "

$c = $ws.Range("D28")
$c.Value = "'161.33"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.81%  "

$c = $ws.Range("D29")
$c.Value = "'2.389"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.00%  "

$c = $ws.Range("D30")
$c.Value = "'129.79"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.65%  "

$c = $ws.Range("D31")
$c.Value = "'1.100"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +3.95%  "

$ws.Range("E32").Value = "  +1.38%  "

$c = $ws.Range("D33")
$c.Value = "'6.048"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.26%  "

$c = $ws.Range("D34")
$c.Value = "'3.730"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.19%  "

$c = $ws.Range("D35")
$c.Value = "'0.02494"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.85%  "

$ws.Range("E36").Value = "  -0.46%  "

$c = $ws.Range("D37")
$c.Value = "'0.2203"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.54%  "

$c = $ws.Range("D38")
$c.Value = "'5.202"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.12%  "

$ws.Range("E39").Value = "  -0.32%  "

$c = $ws.Range("D40")
$c.Value = "'11.91"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +5.72%  "

$c = $ws.Range("D41")
$c.Value = "'8.791"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.12%  "

$c = $ws.Range("D42")
$c.Value = "'0.6526"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.93%  "

$ws.Range("E43").Value = "  -0.24%  "

$c = $ws.Range("D44")
$c.Value = "'0.6124"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.40%  "

$c = $ws.Range("D45")
$c.Value = "'13.18"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.25%  "

$c = $ws.Range("D46")
$c.Value = "'3.720"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.02%  "

$ws.Range("E47").Value = "  +0.99%  "

$c = $ws.Range("D48")
$c.Value = "'1.245"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.32%  "

$c = $ws.Range("D49")
$c.Value = "'124.97"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.51%  "

$c = $ws.Range("D50")
$c.Value = "'1.158"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.82%  "

$c = $ws.Range("D51")
$c.Value = "'79.25"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.14%  "
